# Weekly update: insert 4 new price records (rows 881-884) for "Zapallo" / "Camote"
# at the top of the most-recent-first price history, pushing the existing
# rows 881-929 down to 885-933.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the current row 881. This shifts the
# existing rows 881-929 down to 885-933 (Excel also grows the used range /
# dimension automatically) and the new blank rows inherit the formatting
# (including the date number format on column D) from the row above.
$ws.Rows("881:884").Insert()

# Fill in the 4 new rows with this week's data.
$rows = @(
    @{ Row=881; D=44610; H="Camote"; I="1a (cosecha)";  J=250; K=400; L=430; M=415; O="Región Metropolitana" },
    @{ Row=882; D=44610; H="Camote"; I="1a (cosecha)";  J=160; K=400; L=430; M=415; O="Región de O'Higgins" },
    @{ Row=883; D=44610; H="Camote"; I="2a (cosecha)";  J=133; K=350; L=350; M=350; O="Región Metropolitana" },
    @{ Row=884; D=44610; H="Camote"; I="2a (cosecha)";  J=79;  K=350; L=350; M=350; O="Región de O'Higgins" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112045
    $ws.Cells.Item($row, 7).Value = "Zapallo"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = '$/kilo (volumen en unidades)'
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
